# Commit: "Added dynamic loading calcs"
# Fills in the previously-empty member-force / pin-load columns (F, G, H)
# for the dynamic loading tables (Table 5, rows 66-78, and Table 6, rows
# 85-88) on the "data" worksheet, and updates the saved window/selection
# state to match where the author was last looking in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Table 5: Dynamic Loads - Member forces (rows 66-78) ---------------
$table5 = @{
    66 = @(49513.722756633695,  66587.420258921178,  201469.63052699226)
    67 = @(29708.23365398022,   39952.452155352716,  120881.77831619537)
    68 = @(9902.7445513267448,  13317.484051784249,  40293.926105398481)
    69 = @(-39610.978205306958, -53269.936207136947, -161175.70442159381)
    70 = @(-19805.489102653482, -26634.968103568481, -80587.852210796918)
    71 = @(0, 0, 0)
    72 = @(16504.57425221123,   22195.806752973727,  67156.54350899742)
    73 = @(-19247.475701909207, -25884.536288774452, -78317.314925009865)
    74 = @(19247.475701909207,  25884.536288774452,  78317.314925009865)
    75 = @(-19247.475701909207, -25884.536288774452, -78317.314925009865)
    76 = @(19247.475701909207,  25884.536288774452,  78317.314925009865)
    77 = @(-19247.475701909207, -25884.536288774452, -78317.314925009865)
    78 = @(16504.57425221123,   22195.806752973727,  67156.54350899742)
}

foreach ($row in $table5.Keys) {
    $vals = $table5[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]   # F
    $ws.Cells.Item($row, 7).Value = $vals[1]   # G
    $ws.Cells.Item($row, 8).Value = $vals[2]   # H
}

# --- Table 6: Dynamic Pin loads (rows 85-88) ----------------------------
$table6 = @{
    85 = @(-49513.722756633695, -66587.420258921178, -201469.63052699226)
    86 = @(16504.57425221123,   22195.806752973727,  67156.54350899742)
    87 = @(49513.722756633695,  66587.420258921178,  201469.63052699226)
    88 = @(0, 0, 0)
}

foreach ($row in $table6.Keys) {
    $vals = $table6[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]   # F
    $ws.Cells.Item($row, 7).Value = $vals[1]   # G
    $ws.Cells.Item($row, 8).Value = $vals[2]   # H
}

# --- Update view state: scrolled/selected position on the sheet --------
$ws.Activate()
$ws.Range("L92").Select()
$excel.ActiveWindow.ScrollRow = 66
$excel.ActiveWindow.ScrollColumn = 1

# --- Update the saved workbook window size/position --------------------
$win = $excel.ActiveWindow
$win.WindowState = -4143   # xlNormal
$win.Left = -120
$win.Top = -16320
$win.Width = 29040
$win.Height = 16440
